# Tabela de Testes (VBA) criada - Correção no algoritmo de geração de trajetória
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Start from a clean sheet so stale cells/styles/shared-strings from the
# previous layout (columns B/C/F/G, rows 2-9) are fully removed.
$ws1.Cells.Clear()

# --- Column width / row heights -----------------------------------------
$ws1.Columns.Item(2).ColumnWidth = 26.85546875
$ws1.Rows.Item(2).RowHeight = 18.75
$ws1.Rows.Item(4).RowHeight = 18

# --- Row 1 ----------------------------------------------------------------
$ws1.Range("A1").Value = -5

# --- Row 2 ------------------------------------------------------------------
$ws1.Range("B2").Value = "Initial Acceleration (a0)"
$ws1.Range("B2").Font.Bold = $true
$ws1.Range("C2").Value = 4

$ws1.Range("E2").Value = "(a0)2/(2J)"
$e2 = $ws1.Range("E2")
$e2.Font.Bold = $true
$e2.Characters(3,1).Font.Bold = $true
$e2.Characters(3,1).Font.Subscript = $true
$e2.Characters(4,1).Font.Bold = $true
$e2.Characters(5,1).Font.Bold = $true
$e2.Characters(5,1).Font.Superscript = $true
$e2.Characters(6,5).Font.Bold = $true

$ws1.Range("F2").Value = 1.6

# --- Row 3 ------------------------------------------------------------------
$ws1.Range("B3").Value = "initial Velocity (v0)"
$ws1.Range("B3").Font.Bold = $true
$ws1.Range("C3").Value = 5

$ws1.Range("E3").Value = "v*"
$ws1.Range("E3").Font.Bold = $true

$ws1.Range("F3").Formula = "=C2^2/(2*C5) - C4^2/C5"

# --- Row 4 ------------------------------------------------------------------
$ws1.Range("B4").Value = "Maximum Deceleration (d)"
$ws1.Range("B4").Font.Bold = $true
$ws1.Range("C4").Value = 6

$ws1.Range("E4").Value = "alim"
$e4 = $ws1.Range("E4")
$e4.Font.Bold = $true
$e4.Characters(2,3).Font.Bold = $true
$e4.Characters(2,3).Font.Subscript = $true

# --- Row 5 ------------------------------------------------------------------
$ws1.Range("B5").Value = "Desired Jerk (J)"
$ws1.Range("B5").Font.Bold = $true
$ws1.Range("C5").Value = 5

# --- Row 6 ------------------------------------------------------------------
$ws1.Range("E6").Value = "t1"
$ws1.Range("E6").Font.Bold = $true
$ws1.Range("F6").Formula = "=(C4-C2)/C5"

# --- Row 7 ------------------------------------------------------------------
$ws1.Range("B7").Value = "Initial Position (s0)"
$ws1.Range("B7").Font.Bold = $true
$ws1.Range("C7").Value = 0

$ws1.Range("E7").Value = "t2"
$ws1.Range("E7").Font.Bold = $true
$ws1.Range("F7").Formula = "=C2^2/(2*C5*C4) - C2/C5 - C3/C4"

# --- Row 8 ------------------------------------------------------------------
$ws1.Range("E8").Value = "t3"
$ws1.Range("E8").Font.Bold = $true
$ws1.Range("F8").Formula = "=C4/C5 + F7"

# --- Row 14 -----------------------------------------------------------------
$ws1.Range("J14").NumberFormat = "0.00"

# --- Page setup ---------------------------------------------------------
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# --- Sheet / workbook VBA codenames (best effort; may not persist) -------
try { $wb.CodeName = "EstaPasta_de_trabalho" } catch {}
try { $ws1.CodeName = "Plan1" } catch {}
try { $wb.Worksheets.Item(2).CodeName = "Plan2" } catch {}
try { $wb.Worksheets.Item(3).CodeName = "Plan3" } catch {}
